$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-weight values per "Natmi following Dr Hou advice"
$data = @{
  2 = @{ "E"=2; "G"=0.8352145; "H"=1.670429; "I"=0.07459944913300412; "J"=0.0536415799547362; "K"=2; "M"=64.20135099999999; "N"=128.402702; "O"=0.4070144914449589; "P"=0.3181813759721767; "Q"=53.62189927478949; "R"=214.487597099158; "S"=0.03036305685094376; "T"=0.0170677517193195 }
  3 = @{ "E"=2; "G"=0.8352145; "H"=1.670429; "I"=0.07459944913300412; "J"=0.0536415799547362; "K"=3; "M"=15.75734066666666; "N"=47.27202199999999; "O"=0.09989612209201491; "P"=0.11713987922892; "Q"=13.16075940623966; "R"=78.96455643743799; "S"=0.007452195678587636; "T"=0.006283568197546255 }
  4 = @{ "E"=2; "G"=0.8352145; "H"=1.670429; "I"=0.07459944913300412; "J"=0.0536415799547362; "K"=3; "M"=24.57775566666666; "N"=73.733267; "O"=0.1558145205313015; "P"=0.1827107372630203; "Q"=20.52769791025716; "R"=123.166187461543; "S"=0.01162367739855825; "T"=0.009800892621483101 }
  5 = @{ "E"=2; "G"=0.8352145; "H"=1.670429; "I"=0.07459944913300412; "J"=0.0536415799547362; "K"=3; "M"=14.32600733333333; "N"=42.978022; "O"=0.09082196088386706; "P"=0.1064993646046676; "Q"=11.96528905190633; "R"=71.79173431143799; "S"=0.006775268251115731; "T"=0.005712794181569879 }
  6 = @{ "E"=2; "G"=0.8352145; "H"=1.670429; "I"=0.07459944913300412; "J"=0.0536415799547362; "K"=3; "M"=33.41628466666666; "N"=100.248854; "O"=0.2118477555025799; "P"=0.2484162545532246; "Q"=27.90976548972766; "R"=167.458592938366; "S"=0.0158037258605558; "T"=0.0133254403806729 }
  7 = @{ "E"=2; "G"=0.8352145; "H"=1.670429; "I"=0.07459944913300412; "J"=0.0536415799547362; "K"=2; "M"=5.4585215; "N"=10.917043; "O"=0.03460514954527787; "P"=0.02705238837799083; "Q"=4.55903630536175; "R"=18.236145221447; "S"=0.002581525093242957; "T"=0.001451132854144571 }
  8 = @{ "E"=3; "G"=8.26132; "H"=24.78396; "I"=0.7378822100328354; "J"=0.7958738575150359; "K"=2; "M"=64.20135099999999; "N"=128.402702; "O"=0.4070144914449589; "P"=0.3181813759721767; "Q"=530.3879050433198; "R"=3182.327430259919; "S"=0.3003287524627969; "T"=0.2532322390844182 }
  9 = @{ "E"=3; "G"=8.26132; "H"=24.78396; "I"=0.7378822100328354; "J"=0.7958738575150359; "K"=3; "M"=15.75734066666666; "N"=47.27202199999999; "O"=0.09989612209201491; "P"=0.11713987922892; "Q"=130.1764335963466; "R"=1171.58790236712; "S"=0.07371157134296591; "T"=0.093228567550766 }
  10 = @{ "E"=3; "G"=8.26132; "H"=24.78396; "I"=0.7378822100328354; "J"=0.7958738575150359; "K"=3; "M"=24.57775566666666; "N"=73.733267; "O"=0.1558145205313015; "P"=0.1827107372630203; "Q"=203.0447044441466; "R"=1827.40233999732; "S"=0.1149727627648433; "T"=0.1454146992749362 }
  11 = @{ "E"=3; "G"=8.26132; "H"=24.78396; "I"=0.7378822100328354; "J"=0.7958738575150359; "K"=3; "M"=14.32600733333333; "N"=42.978022; "O"=0.09082196088386706; "P"=0.1064993646046676; "Q"=118.3517309030133; "R"=1065.16557812712; "S"=0.06701590921650355; "T"=0.08476006013081708 }
  12 = @{ "E"=3; "G"=8.26132; "H"=24.78396; "I"=0.7378822100328354; "J"=0.7958738575150359; "K"=3; "M"=33.41628466666666; "N"=100.248854; "O"=0.2118477555025799; "P"=0.2484162545532246; "Q"=276.0626208424266; "R"=2484.56358758184; "S"=0.1563186900207394; "T"=0.197708002780712 }
  13 = @{ "E"=3; "G"=8.26132; "H"=24.78396; "I"=0.7378822100328354; "J"=0.7958738575150359; "K"=2; "M"=5.4585215; "N"=10.917043; "O"=0.03460514954527787; "P"=0.02705238837799083; "Q"=45.09459283838; "R"=270.56755703028; "S"=0.0255345242249864; "T"=0.02153028869338649 }
  14 = @{ "E"=3; "G"=0.487268; "H"=1.461804; "I"=0.04352166345308978; "J"=0.04694211854808148; "K"=2; "M"=64.20135099999999; "N"=128.402702; "O"=0.4070144914449589; "P"=0.3181813759721767; "Q"=31.28326389906799; "R"=187.6995833944079; "S"=0.01771394771719799; "T"=0.0149361078706776 }
  15 = @{ "E"=3; "G"=0.487268; "H"=1.461804; "I"=0.04352166345308978; "J"=0.04694211854808148; "K"=3; "M"=15.75734066666666; "N"=47.27202199999999; "O"=0.09989612209201491; "P"=0.11713987922892; "Q"=7.678047871965332; "R"=69.10243084768798; "S"=0.004347645405957439; "T"=0.00549879409747191 }
  16 = @{ "E"=3; "G"=0.487268; "H"=1.461804; "I"=0.04352166345308978; "J"=0.04694211854808148; "K"=3; "M"=24.57775566666666; "N"=73.733267; "O"=0.1558145205313015; "P"=0.1827107372630203; "Q"=11.97595384818533; "R"=107.783584633668; "S"=0.00678130712366785; "T"=0.008576829088608066 }
  17 = @{ "E"=3; "G"=0.487268; "H"=1.461804; "I"=0.04352166345308978; "J"=0.04694211854808148; "K"=3; "M"=14.32600733333333; "N"=42.978022; "O"=0.09082196088386706; "P"=0.1064993646046676; "Q"=6.980604941298665; "R"=62.82544447168799; "S"=0.003952722815737346; "T"=0.004999305798567659 }
  18 = @{ "E"=3; "G"=0.487268; "H"=1.461804; "I"=0.04352166345308978; "J"=0.04694211854808148; "K"=3; "M"=33.41628466666666; "N"=100.248854; "O"=0.2118477555025799; "P"=0.2484162545532246; "Q"=16.28268619695733; "R"=146.544175772616; "S"=0.00921996671827573; "T"=0.01166118527050786 }
  19 = @{ "E"=3; "G"=0.487268; "H"=1.461804; "I"=0.04352166345308978; "J"=0.04694211854808148; "K"=2; "M"=5.4585215; "N"=10.917043; "O"=0.03460514954527787; "P"=0.02705238837799083; "Q"=2.659762854262; "R"=15.958577125572; "S"=0.001506073672253426; "T"=0.001269896422248387 }
  20 = @{ "E"=2; "G"=1.612185; "H"=3.22437; "I"=0.1439966773810707; "J"=0.1035424439821464; "K"=2; "M"=64.20135099999999; "N"=128.402702; "O"=0.4070144914449589; "P"=0.3181813759721767; "Q"=103.504455061935; "R"=414.0178202477399; "S"=0.05860873441402031; "T"=0.03294527729776136 }
  21 = @{ "E"=2; "G"=1.612185; "H"=3.22437; "I"=0.1439966773810707; "J"=0.1035424439821464; "K"=3; "M"=15.75734066666666; "N"=47.27202199999999; "O"=0.09989612209201491; "P"=0.11713987922892; "Q"=25.40374826268999; "R"=152.42248957614; "S"=0.01438470966450392; "T"=0.01212894938313584 }
  22 = @{ "E"=2; "G"=1.612185; "H"=3.22437; "I"=0.1439966773810707; "J"=0.1035424439821464; "K"=3; "M"=24.57775566666666; "N"=73.733267; "O"=0.1558145205313015; "P"=0.1827107372630203; "Q"=39.623889019465; "R"=237.74333411679; "S"=0.02243677324423203; "T"=0.01891831627799294 }
  23 = @{ "E"=2; "G"=1.612185; "H"=3.22437; "I"=0.1439966773810707; "J"=0.1035424439821464; "K"=3; "M"=14.32600733333333; "N"=42.978022; "O"=0.09082196088386706; "P"=0.1064993646046676; "Q"=23.09617413269; "R"=138.57704479614; "S"=0.01307806060051043; "T"=0.01102720449371298 }
  24 = @{ "E"=2; "G"=1.612185; "H"=3.22437; "I"=0.1439966773810707; "J"=0.1035424439821464; "K"=3; "M"=33.41628466666666; "N"=100.248854; "O"=0.2118477555025799; "P"=0.2484162545532246; "Q"=53.87323289532999; "R"=323.23939737198; "S"=0.03050537290300894; "T"=0.02572162612133188 }
  25 = @{ "E"=2; "G"=1.612185; "H"=3.22437; "I"=0.1439966773810707; "J"=0.1035424439821464; "K"=2; "M"=5.4585215; "N"=10.917043; "O"=0.03460514954527787; "P"=0.02705238837799083; "Q"=8.8001464844775; "R"=35.20058593791; "S"=0.004983026554795082; "T"=0.002801070408211383 }
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  foreach ($col in $row.Keys) {
    $ws.Range("$col$r").Value = $row[$col]
  }
}
